# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.916.70"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.660.08"
$ws.Range("E3").Value = "  +2.12%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("E5").Value = "  +1.41%  "
$ws.Range("D6").Value = "'0.520"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'28.94"
$ws.Range("E8").Value = "  -1.98%  "
$ws.Range("D9").Value = "'0.263"
$ws.Range("E9").Value = "  +2.01%  "
$ws.Range("D10").Value = "'0.0614"
$ws.Range("E10").Value = "  +0.45%  "
$ws.Range("D11").Value = "'0.0901"
$ws.Range("E11").Value = "  -1.43%  "
$ws.Range("D12").Value = "1.896.17"
$ws.Range("E12").Value = "  +2.00%  "
$ws.Range("D13").Value = "1.658.17"
$ws.Range("E13").Value = "  +1.86%  "
$ws.Range("D14").Value = "'0.601"
$ws.Range("E14").Value = "  +4.97%  "
$ws.Range("D15").Value = "'10.04"
$ws.Range("E15").Value = "  +12.90%  "
$ws.Range("D16").Value = "'3.95"
$ws.Range("E16").Value = "  +1.31%  "
$ws.Range("D17").Value = "29.922.92"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "'64.72"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").Value = "'241.53"
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").Value = "'10.02"
$ws.Range("E22").Value = "  +4.07%  "
$ws.Range("D23").Value = "'4.17"
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("D24").Value = "'2.19"
$ws.Range("E24").Value = "  +2.93%  "
$ws.Range("D25").Value = "'158.05"
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("D26").Value = "'15.73"
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("E30").Value = "  +1.77%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").Value = "'3.40"
$ws.Range("E32").Value = "  +1.69%  "
$ws.Range("D33").Value = "'3.22"
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("D34").Value = "1.438.90"
$ws.Range("E34").Value = "  +1.22%  "
$ws.Range("D35").Value = "'1.71"
$ws.Range("E35").Value = "  +4.55%  "
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("D37").Value = "'0.0175"
$ws.Range("E37").Value = "  +2.59%  "
$ws.Range("B38").Value = "Aave"
$ws.Range("C38").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D38").Value = "'78.45"
$ws.Range("E38").Value = "  +13.12%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.68"
$ws.Range("E39").Value = "  -6.53%  "
$ws.Range("D40").Value = "'0.574"
$ws.Range("E40").Value = "  +3.00%  "
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").Value = "'0.844"
$ws.Range("E42").Value = "  +1.52%  "
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("E44").Value = "  -2.13%  "
$ws.Range("D45").Value = "'0.998"
$ws.Range("E46").Value = "  -1.94%  "
$ws.Range("D47").Value = "'50.75"
$ws.Range("E47").Value = "  -6.42%  "
$ws.Range("D48").Value = "1.803.17"
$ws.Range("E48").Value = "  +2.05%  "
$ws.Range("D49").Value = "'5.37"
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").Value = "'94.20"
$ws.Range("E50").Value = "  +6.17%  "
$ws.Range("D51").Value = "0.0₆0109"
$ws.Range("E51").Value = "  +0.66%  "
